# Update computed market-price / profit figures on each crafting-Leve sheet.
# Values sourced from a scheduled Universalis price-refresh run (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3: One for the Books
$ws.Range("H3").Value = 29800
$ws.Range("J3").Value = 29800
$ws.Range("L3").Value = 29800
$ws.Range("N3").Value = -30028

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 1392.4546
$ws.Range("I40").Value = 650
$ws.Range("J40").Value = 1557.4445
$ws.Range("K40").Value = 650
$ws.Range("L40").Value = 1557.4445
$ws.Range("M40").Value = -475
$ws.Range("N40").Value = -1907.4445

# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 411.6111
$ws.Range("I92").Value = 394.3125
$ws.Range("K92").Value = 394.3125
$ws.Range("M92").Value = 853.6875

# Row 102: Spell-rebound
$ws.Range("H102").Value = 29800
$ws.Range("J102").Value = 29800
$ws.Range("L102").Value = 29800
$ws.Range("N102").Value = -36290

# Row 109: A Time for Peace
$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774

# Row 116: Growing Up
$ws.Range("H116").Value = 3937.8125
$ws.Range("I116").Value = 2166.6667
$ws.Range("K116").Value = 2166.6667
$ws.Range("M116").Value = 1275.3333

# Row 123: Nearly Bare
$ws.Range("H123").Value = 30780
$ws.Range("J123").Value = 30780
$ws.Range("L123").Value = 30780
$ws.Range("N123").Value = -40580

# Row 129: Practical Command
$ws.Range("H129").Value = 756.5946
$ws.Range("I129").Value = 375
$ws.Range("J129").Value = 802.8484999999999
$ws.Range("K129").Value = 1125
$ws.Range("L129").Value = 2408.5455
$ws.Range("M129").Value = 3875
$ws.Range("N129").Value = -12408.5455

# Row 135: For Tired Minds
$ws.Range("H135").Value = 20838590
$ws.Range("I135").Value = 1425.4166
$ws.Range("J135").Value = 41675756
$ws.Range("K135").Value = 12828.7494
$ws.Range("L135").Value = 375081804
$ws.Range("M135").Value = -10293.7494
$ws.Range("N135").Value = -375086874

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1546.2963
$ws.Range("I137").Value = 1663.8889
$ws.Range("J137").Value = 1311.1111
$ws.Range("K137").Value = 4991.6667
$ws.Range("L137").Value = 3933.3333
$ws.Range("M137").Value = -2441.6667
$ws.Range("N137").Value = -9033.3333

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2173.1592
$ws.Range("J138").Value = 3386.95
$ws.Range("L138").Value = 10160.85
$ws.Range("N138").Value = -20440.85

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7436.2837
$ws.Range("I32").Value = 6177.3
$ws.Range("K32").Value = 6177.3
$ws.Range("M32").Value = -5890.3

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 3094.6
$ws.Range("I45").Value = 2590.375
$ws.Range("J45").Value = 3670.8572
$ws.Range("K45").Value = 2590.375
$ws.Range("L45").Value = 3670.8572
$ws.Range("M45").Value = -2213.375
$ws.Range("N45").Value = -4424.8572

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 23811870
$ws.Range("I74").Value = 34484708
$ws.Range("K74").Value = 34484708
$ws.Range("M74").Value = -34483834

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 23811870
$ws.Range("I77").Value = 34484708
$ws.Range("K77").Value = 172423540
$ws.Range("M77").Value = -172419172

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1423.1333
$ws.Range("I102").Value = 1364.3077
$ws.Range("K102").Value = 1364.3077
$ws.Range("M102").Value = 257.6922999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1831.3462
$ws.Range("I86").Value = 1652.3334
$ws.Range("K86").Value = 1652.3334
$ws.Range("M86").Value = -529.3334

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1831.3462
$ws.Range("I89").Value = 1652.3334
$ws.Range("K89").Value = 8261.666999999999
$ws.Range("M89").Value = -2645.666999999999

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 1614721
$ws.Range("I105").Value = 1654.1666
$ws.Range("K105").Value = 1654.1666
$ws.Range("M105").Value = 92.83339999999998

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 3898.25
$ws.Range("I134").Value = 3571.4707
$ws.Range("K134").Value = 10714.4121
$ws.Range("M134").Value = -8179.4121

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 4153.7085
$ws.Range("I31").Value = 799.4
$ws.Range("J31").Value = 6549.643
$ws.Range("K31").Value = 799.4
$ws.Range("L31").Value = 6549.643
$ws.Range("M31").Value = -504.4
$ws.Range("N31").Value = -7139.643

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 4153.7085
$ws.Range("I34").Value = 799.4
$ws.Range("J34").Value = 6549.643
$ws.Range("K34").Value = 799.4
$ws.Range("L34").Value = 6549.643
$ws.Range("M34").Value = -597.4
$ws.Range("N34").Value = -6953.643

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1176.2354
$ws.Range("I105").Value = 988.3077
$ws.Range("J105").Value = 1787
$ws.Range("K105").Value = 988.3077
$ws.Range("L105").Value = 1787
$ws.Range("M105").Value = 758.6923
$ws.Range("N105").Value = -5281

# Row 107: Built to Last
$ws.Range("H107").Value = 1129.1538
$ws.Range("I107").Value = 412.35715
$ws.Range("K107").Value = 412.35715
$ws.Range("M107").Value = 1507.64285

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch
$ws.Range("H34").Value = 543.3
$ws.Range("I34").Value = 105.8
$ws.Range("J34").Value = 980.8
$ws.Range("K34").Value = 317.4
$ws.Range("L34").Value = 2942.4
$ws.Range("M34").Value = -233.4
$ws.Range("N34").Value = -3110.4

# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 3337.9
$ws.Range("J39").Value = 3337.9
$ws.Range("L39").Value = 10013.7
$ws.Range("N39").Value = -10601.7

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 2639.5
$ws.Range("J55").Value = 3276.875
$ws.Range("L55").Value = 9830.625
$ws.Range("N55").Value = -10184.625

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 796.2222
$ws.Range("I113").Value = 616
$ws.Range("J113").Value = 1156.6666
$ws.Range("K113").Value = 1848
$ws.Range("L113").Value = 3469.9998
$ws.Range("M113").Value = 322
$ws.Range("N113").Value = -7809.9998

# Row 122: Salt of the North
$ws.Range("H122").Value = 1150
$ws.Range("I122").Value = 207.4
$ws.Range("J122").Value = 1398.0526
$ws.Range("K122").Value = 1866.6
$ws.Range("L122").Value = 12582.4734
$ws.Range("M122").Value = 583.3999999999999
$ws.Range("N122").Value = -17482.4734

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 725.71
$ws.Range("J131").Value = 730.0103
$ws.Range("L131").Value = 2190.0309
$ws.Range("N131").Value = -12270.0309

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 457.375
$ws.Range("J107").Value = 621.6667
$ws.Range("L107").Value = 621.6667
$ws.Range("N107").Value = -4461.6667

# Row 132: On Board for Lar
$ws.Range("H132").Value = 21732.178
$ws.Range("I132").Value = 3985.3635
$ws.Range("J132").Value = 86803.836
$ws.Range("K132").Value = 11956.0905
$ws.Range("L132").Value = 260411.508
$ws.Range("M132").Value = -9426.0905
$ws.Range("N132").Value = -265471.508

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 4582.231
$ws.Range("I7").Value = 4187
$ws.Range("K7").Value = 4187
$ws.Range("M7").Value = -4075

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 3224.9167
$ws.Range("I22").Value = 3633.4443
$ws.Range("K22").Value = 3633.4443
$ws.Range("M22").Value = -3338.4443

# Row 27: Fire and Hide
$ws.Range("H27").Value = 3224.9167
$ws.Range("I27").Value = 3633.4443
$ws.Range("K27").Value = 3633.4443
$ws.Range("M27").Value = -3526.4443

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 1102.6459
$ws.Range("I46").Value = 1085.909
$ws.Range("J46").Value = 1286.75
$ws.Range("K46").Value = 1085.909
$ws.Range("L46").Value = 1286.75
$ws.Range("M46").Value = -897.9090000000001
$ws.Range("N46").Value = -1662.75

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 3477.9546
$ws.Range("I61").Value = 1362
$ws.Range("J61").Value = 12999.75
$ws.Range("K61").Value = 1362
$ws.Range("L61").Value = 12999.75
$ws.Range("M61").Value = -1160
$ws.Range("N61").Value = -13403.75

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 2249.5
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 2666
$ws.Range("K68").Value = 1000
$ws.Range("L68").Value = 2666
$ws.Range("M68").Value = -251
$ws.Range("N68").Value = -4164

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 2249.5
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 2666
$ws.Range("K71").Value = 5000
$ws.Range("L71").Value = 13330
$ws.Range("M71").Value = -1256
$ws.Range("N71").Value = -20818

# Row 113: Peace in Rest
$ws.Range("H113").Value = 3477.9546
$ws.Range("I113").Value = 1362
$ws.Range("J113").Value = 12999.75
$ws.Range("K113").Value = 1362
$ws.Range("L113").Value = 12999.75
$ws.Range("M113").Value = 808
$ws.Range("N113").Value = -17339.75

# Row 122: Hell on Leather
$ws.Range("H122").Value = 656607.0600000001
$ws.Range("I122").Value = 728907.9
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 2186723.7
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -2184273.7
$ws.Range("N122").Value = -22600

# Row 126: Battered Books
$ws.Range("H126").Value = 4582.231
$ws.Range("I126").Value = 4187
$ws.Range("K126").Value = 12561
$ws.Range("M126").Value = -10091

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 1535.037
$ws.Range("I136").Value = 1363.3077
$ws.Range("K136").Value = 4089.9231
$ws.Range("M136").Value = -1539.9231

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 2262.9285
$ws.Range("I81").Value = 1566.6666
$ws.Range("J81").Value = 2785.125
$ws.Range("K81").Value = 3133.3332
$ws.Range("L81").Value = 5570.25
$ws.Range("M81").Value = -2072.3332
$ws.Range("N81").Value = -7692.25

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 2262.9285
$ws.Range("I84").Value = 1566.6666
$ws.Range("J84").Value = 2785.125
$ws.Range("K84").Value = 15666.666
$ws.Range("L84").Value = 27851.25
$ws.Range("M84").Value = -10362.666
$ws.Range("N84").Value = -38459.25

# Row 100: Of Great Import
$ws.Range("H100").Value = 219.3077
$ws.Range("I100").Value = 194.9
$ws.Range("K100").Value = 389.8
$ws.Range("M100").Value = 151.2

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1157.75
$ws.Range("I132").Value = 952.7917
$ws.Range("J132").Value = 1465.1875
$ws.Range("K132").Value = 2858.3751
$ws.Range("L132").Value = 4395.5625
$ws.Range("M132").Value = -328.3751000000002
$ws.Range("N132").Value = -9455.5625
